$wb = $excel.ActiveWorkbook
$sheet4 = $wb.Worksheets.Item("Remediation")

# --- Update existing rows 2-4, column A: change from numeric id to a text bug id ---
# Order chosen to reproduce the exact shared-string table ordering of the target file.
$sheet4.Range("A8").Value = "json"
$sheet4.Range("C5").Value = "issues with calls to the mongo db not returning valid json"
$sheet4.Range("B5").Value = "json1"
$sheet4.Range("A3").Value = "ui2"
$sheet4.Range("A4").Value = "ui3"
$sheet4.Range("A2").Value = "field1"
$sheet4.Range("A6").Value = "ui1"
$sheet4.Range("B6").Value = "login"
$sheet4.Range("C6").Value = "login via spotify does not display welcome message in correct place"
$sheet4.Range("E5").Value = "Joe/Cullen"
$sheet4.Range("A7").Value = "ui4"
$sheet4.Range("B7").Value = "rsearches"
$sheet4.Range("C7").Value = "recent search results not displaying in columns"
$sheet4.Range("A5").Value = "mongo1"
$sheet4.Range("B8").Value = "json2"
$sheet4.Range("C8").Value = 'index/home not property passing "similar tracks" to similar tracks page'

# --- Remaining cells of the 4 new rows (5-8) ---
$sheet4.Range("D5").Value = 3
$sheet4.Range("D6").Value = 2
$sheet4.Range("D7").Value = 2
$sheet4.Range("D8").Value = 3

$sheet4.Range("E6").Value = "Joe"
$sheet4.Range("E7").Value = "Joe"
$sheet4.Range("E8").Value = "Cullen"

# Dates in column F re-use the existing date style (numFmtId 14) already present on F2:F4,
# so copy that formatting across instead of assigning a brand new number format.
$sheet4.Range("F2").Copy() | Out-Null
$sheet4.Range("F5:F8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$sheet4.Range("F5").Value = 42709
$sheet4.Range("F6").Value = 42709
$sheet4.Range("F7").Value = 42709
$sheet4.Range("F8").Value = 42709

# --- Restore the active selection shown in the saved workbook ---
$sheet4.Range("F12").Select()
